$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$A203 = @'
Right. So absolutely you know but in the past you know I would say since Air India 171 crash uh there are many developments that have taken place uh uh uh you know whether it's about this one or the ones in the past. Uh many many incidents are being given limelight. Many many incidents are coming to the front. uh minor ones of course but but you know given the situation right now uh what exactly is your assessment on that front sir I think this was the right situation because uh I I stay in Bombay and I've seen this morning there were very heavy rains water logging all around so in my opinion there must have been huge quantity of run of the water on the runway and I really don't know uh I'm just hoping like hell that it was safe safer level of water over the runway for landing uh and why pilot uh whether pilot thought of diverting uh to the the destination or not. Uh all this thing will come out in the in the investigation report. All right. Absolutely. Uh with that, I'd like to thank you for joining us on the broadcast, sharing your perspectives with us on NewsX. For more such videos, subscribe to the NewsX YouTube channel. Hit the bell icon. [Music]
'@
$B203 = @'
since air india 171 crash there have been many developments that have taken place many many incidents are coming to the front i think this was the right situation because there were very heavy rains water logging all around all this thing will come out in the in the investigation report
'@
$C203 = @'
political
'@
$A204 = @'
China has begun construction on one of its most ambitious infrastructure projects ever. A $170 billion migdum on the Yalong Zangbo River just upstream of India's Arunachal Pradesh. This dam spread across a 50 km stretch will eventually generate 60 gawatt of electricity nearly three times the capacity of the iconic three gorgeous dam. But there's more than just energy at stake. The Yalang Zangbo becomes the Brahmaputra as it enters India through Arunachal. And right at this critical bend, Beijing is building a series of five casket dams with little information shared about their water storage capacity. something that directly affects India's river flow. For New Delhi, the timing and scale of this announcement are impossible to ignore. Back home, alarm bells are ringing, especially in Arunachal Pradesh. Chief Minister Parandu has called the Chinese project a bigger threat than even military incursions. His warning is stuck. If China suddenly releases water, the entire sigh belt in Arunachal, home to indigenous communities like the Adi tribe could be wiped out. Floods, land loss, and human displacement, are no longer distant fears. They're emerging realities. The stakes go beyond floods. Arunachal alone holds nearly 50 gawatt of untapped hydropower potential. That's over 80% of what the entire northeast offers. China's upstream control could directly impact India's ability to harness this power, disrupting flow, delaying projects, and undermining long-term water security in the region. To counter Beijing's water leverage, India has a plan of its own. The Upper Sai multi-purpose project, a massive 11.2 gawatt dam in Arunachal, is designed to regulate river flow and serve as a strategic buffer. It's intended to act as both a hydropower asset and a safeguard against floods or artificial water releases by China. But there's a catch. The project has seen little progress. Prefeasibility studies are stalled due to local opposition and ground level investigations haven't taken off. Despite the delay, the center says it's committed. Union Jal Shakti Minister CR Patel told the Indian Express earlier this month, "China can do whatever it wants. We are fully prepared." PM Modi is quite serious about it. Work will start. The NHPC, India's hydropower major, maintains that the upper dam can act as a regulating mechanism, potentially blunting any water shocks from Chinese infrastructure upstream. As China powers ahead with its $170 billion gamble, India faces a strategic test. Balancing local resistance, environmental sensitivities and national security concerns in the race to secure river control, water may soon become as critical a frontier as land or air. Whether India can act before the floodgates open, that remains to be seen. Heat. Heat. [Music]
'@
$B204 = @'
china has begun construction on one of its most ambitious infrastructure projects ever a 170 billion migdum on the yalong zangbo river just upstream of indias arunachal pradesh this dam spread across a 50 km stretch will eventually generate 60 gawatt of electricity
'@
$C204 = @'
political
'@
$A205 = @'
Shortly after takeoff, a Bangladesh Air Force training aircraft crashed into a school campus in the capital, Dhaka on July 21st. 19 people, including 16 students, two teachers and a pilot have been killed in the crash. The China made fighter jet crashed into the campus of Milestone School and College in Dhaka, severely injuring more than 92 people. Eight still remain critical. They have been rushed to nearby medical facilities and are being treated for severe burns. The country's interim leader Muhammad Yunus expressed deep grief and sorrow over the incident in a post on X. According to the Associated Press, the jet took off at 10:06 p.m. and crashed into the school packed with students shortly after it caught fire. The jet was an F7BGI, an advanced version of the Chinese J7 fighter and Bangladesh Air Force has 16 of these. This is said to be the deadliest airplane crash in the Bangladeshi capital in recent memory. Panicked students, parents, and teachers were seen running in confusion with a backdrop of screams, whales, and raging orange flames. The student's body was seen lying in a pool of blood as a policeman rushed towards the wreckage, apparently searching for survivors in one of several videos emerging from the tragedy. Authorities fear that the death toll may rise. Bodies of those who were killed in the crash will be handed over to their families immediately. [Music]
'@
$B205 = @'
19 people including 16 students two teachers and a pilot have been killed in the crash the china made fighter jet crashed into the campus of milestone school and college in dhaka severely injuring more than 92 people eight still remain critical this is said to be the deadliest airplane crash in the bangladeshi capital in recent memory
'@
$C205 = @'
technology 
'@

$ws.Range("A203").Value = $A203
$ws.Range("B203").Value = $B203
$ws.Range("C203").Value = $C203

$ws.Range("A204").Value = $A204
$ws.Range("B204").Value = $B204
$ws.Range("C204").Value = $C204

$ws.Range("A205").Value = $A205
$ws.Range("B205").Value = $B205
$ws.Range("C205").Value = $C205
